$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.642.49"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.099.58"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.19"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.55"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.095.54"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.435"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.109"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("D13").Value = "3.629.52"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.27"
$ws.Range("E15").Value = "  +3.59%  "
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "57.735.09"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "3.097.73"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.81"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.06"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "337.26"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.506"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.11"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "0.0₃0915"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.57"
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.95"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  +3.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.06"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0664"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "3.140.66"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.89"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.78"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  +6.66%  "
$ws.Range("D47").Value = "2.280.61"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.60"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.959"
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.01"
$ws.Range("E51").Value = "  +2.65%  "
